$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "vahan.hovhannisyan.h@gmail.com" (row 2),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Select cell B6 as the active cell, matching the saved view state.
$ws.Range("B6").Select()
